$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.786.47"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "3.885.10"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.81"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.53"
$ws.Range("E6").Value = "  +9.74%  "
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.760"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("E10").Value = "  +6.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.01"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.33"
$ws.Range("E13").Value = "  +5.50%  "
$ws.Range("D14").Value = "4.512.93"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").Value = "3.883.47"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.00"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").Value = "71.660.88"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.21"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.19"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("E24").Value = "  -4.09%  "
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.18"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  -4.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.08"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("E31").Value = "  +4.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.21"
$ws.Range("E32").Value = "  +5.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "13.58"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("E34").Value = "  -3.46%  "
$ws.Range("D35").Value = "0.0₃0979"
$ws.Range("E35").Value = "  +14.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "68.37"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "622.82"
$ws.Range("E37").Value = "  -8.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.420"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.30"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.22"
$ws.Range("E43").Value = "  +41.29%  "
$ws.Range("E44").Value = "  -2.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.25"
$ws.Range("E45").Value = "  -6.15%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.144"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.63"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.82"
$ws.Range("E49").Value = "  -15.51%  "
$ws.Range("D50").Value = "2.870.87"
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("E51").Value = "  +1.10%  "
